# "Premier test guide signalisation"
# Update signage table on sheet "Etape1" (physically sheet3.xml) with the
# revised / corrected entries for rows 5-24 (columns D, E, F, G), then move
# the active selection to F24 as in the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Etape1")

# --- Column D (details) updates (ordered to mirror the authored commit's
#     shared-string insertion order) ------------------------------------
$ws.Range("D5").Value  = "Contrôle policier, intersection ch mine Westwood"
$ws.Range("D6").Value  = "Contrôle policier, intersection chemin Preissac"
$ws.Range("D7").Value  = "Contrôle policier, intersection secteur Cadillac (École)"
$ws.Range("D8").Value  = "Contrôle policier, intersection secteur Cadillac (Station-service)"
$ws.Range("D13").Value = "Terre-Plein entrée rond point Malartic - tous à droite <br/>Série de cônes vers la direction droite"
$ws.Range("D15").Value = "Voie ferrée - Malartic"
$ws.Range("D19").Value = "Retressissment Ville de Malartic (Caisse)"
$ws.Range("D16").Value = "Retressissment Ville de Malartic (Pétro)"
$ws.Range("D23").Value = "Retressissment Ville de Malartic - terre-plein sortie ville"
$ws.Range("D22").Value = "Retressissment Ville de Malartic - terre-plein av Hochelaga"

# --- Column E (fonction) updates ---------------------------------------
$ws.Range("E11").Value = "signaleur_moto"
$ws.Range("E13").Value = "signaleur_moto"

# --- Column F (type) updates --------------------------------------------
$ws.Range("F13").Value = "terre_plein"
$ws.Range("F24").Value = "terre_plein"

# --- Column G (responsable) updates -------------------------------------
$ws.Range("G9").Value  = "sq_autre"
$ws.Range("G10").Value = "sq_autre"
$ws.Range("G11").Value = "signaleur_moto"
$ws.Range("G12").Value = "sq_hotesse"
$ws.Range("G13").Value = "signaleur_moto"
$ws.Range("G14").Value = "sq_hotesse"

# --- Update active selection to F24, matching the authored edit ---------
$ws.Activate()
$ws.Range("F24").Select()
